# Ajustes a flujo de recordar usuario
#
# Updates the "Datos" sheet of the recordar_usuario data-driven test
# workbook: new numeroDocumento / usuario test values per row, a
# corrected clave on row 3, idCaso on row 4 stored as text, shorter /
# swapped mensajeRespuesta texts, row 6 losing its stray highlight
# style, and the active selection moving to K12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Note: the K column (and A4 below) already carries a quoted-text style
# (quotePrefix) in the template. A leading apostrophe keeps Excel's
# "stored as text" marker on the cell (same style) while still writing
# the clean text as the value, just like typing `'some text` by hand.

# --- Row 2 (idCaso 1 - Acierto) ---------------------------------------
$ws.Range("F2").Value = "25130110"
$ws.Range("G2").Value = "USUCTDC1"
$ws.Range("K2").Value = "'El usuario ha sido enviado al correo electrónico."

# --- Row 3 (idCaso 2 - Alterno) ---------------------------------------
$ws.Range("F3").Value = "25130112"
$ws.Range("G3").Value = "USUCTDC3"
$ws.Range("H3").Value = 4321
$ws.Range("K3").Value = "'Usuario o clave inválida. Inténtalo nuevamente"

# --- Row 4 (idCaso 3 - Alterno) ---------------------------------------
# idCaso becomes a text value instead of a number.
$ws.Range("A4").Value = "'3"
$ws.Range("F4").Value = "25130111"
$ws.Range("G4").Value = "USUCTDC2"
$ws.Range("K4").Value = "'¡Lo Sentimos!"

# --- Row 5 (idCaso 4 - Alterno) ---------------------------------------
$ws.Range("F5").Value = "1989636240"
$ws.Range("G5").Value = "OSVPPRU16"
$ws.Range("K5").Value = "'La clave que usas en el cajero está bloqueada."

# --- Row 6 (idCaso 5 - Alterno) ---------------------------------------
# F6/G6 drop the extra highlighted style and line up with the other
# rows' numeroDocumento/usuario formatting.
$ws.Range("F6").NumberFormat = $ws.Range("F2").NumberFormat
$ws.Range("G6").NumberFormat = $ws.Range("G2").NumberFormat
$ws.Range("F6").Value = "25130114"
$ws.Range("G6").Value = "USUCTDC5"

# --- Active selection ---------------------------------------------------
$ws.Range("K12").Select()
